$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update numeric result cells (column C) and vector/result strings ---
$ws.Range("C3").Value = 2
$ws.Range("F3").Value = '14.9847           15      3.01046'
$ws.Range("C4").Value = 1.8973
$ws.Range("F4").Value = '7.84644      4.58378      12.9908'
$ws.Range("C5").Value = 1.5851
$ws.Range("F5").Value = '3.2827        6.95      6.0672'
$ws.Range("C7").Value = 1.4815499999999999
$ws.Range("D7").Value = '0.751498     0.676485     0.878121 '
$ws.Range("F7").Value = '13.0375      14.9529      3.58575'
$ws.Range("C8").Value = 1.3889
$ws.Range("D8").Value = '0.74946     0.79341     0.82094 '
$ws.Range("F8").Value = ' 6.846      2.8914      7.0493'
$ws.Range("C9").Value = 1.7917000000000001
$ws.Range("D9").Value = '0.9958    0.6091    0.8745'
$ws.Range("F9").Value = '11.9855    4.2149    4.5718'
$ws.Range("C11").Value = 1.5340499999999999
$ws.Range("D11").Value = '0.75196     0.564448     0.960431'
$ws.Range("F11").Value = '9.84823      18.5548      5.14193 '
$ws.Range("I11").Value = ' 60.5674      60.4803      59.9803'
$ws.Range("J11").Value = '89.6256      88.6417      91.7912'
$ws.Range("C12").Value = 1.9049739999999999
$ws.Range("D12").Value = '0.9160724     0.8638073     0.5727551  '
$ws.Range("F12").Value = '10.69823      3.682638      7.553252'
$ws.Range("I12").Value = '53.38334      40.00147      51.29529 '
$ws.Range("J12").Value = '96.4811      92.42004      107.2134'
$ws.Range("C13").Value = 1.876174
$ws.Range("D13").Value = '0.7378902     0.7642485     0.7004045'
$ws.Range("F13").Value = '17.6893       15.8735      13.69968 '
$ws.Range("I13").Value = '59.56614      66.37721      46.04473 '
$ws.Range("J13").Value = '108.7981      131.4733      94.46455'
$ws.Range("C15").Value = 1.92988
$ws.Range("D15").Value = '0.74941      0.82108     0.733546  '
$ws.Range("F15").Value = '7.5277      16.5885      8.62631'
$ws.Range("C16").Value = 1.6651400000000001
$ws.Range("D16").Value = ' 0.727287     0.730442     0.721689 '
$ws.Range("F16").Value = '6.1908      8.55466      15.6824'
$ws.Range("L16").Value = '0, 0.025, 0.045'
$ws.Range("C17").Value = 1.8469
$ws.Range("D17").Value = ' 0.73358     0.83923     0.78457'
$ws.Range("F17").Value = '0.10025      2.8951      2.4026'
$ws.Range("C18").Value = 1.9643999999999999
$ws.Range("D18").Value = ' 0.75529           1      0.5876'
$ws.Range("F18").Value = '1.066      2.1154      1.3248'
$ws.Range("C20").Value = '0.45319      1.3363      2.8087 '
$ws.Range("F20").Value = '5.0797      2.8558      1.5681'
$ws.Range("G20").Value = '5.0797      1.6178      6.7134'
$ws.Range("H20").Value = '5.0797      5.2114      4.2905'

# --- Update the active cell selection to match the post-edit state ---
$ws.Range("H20").Select()
